$wb = $excel.ActiveWorkbook

# Sheet "navi"
$ws1 = $wb.Worksheets.Item("navi")
$ws1.Range("C2").Value = 30
$ws1.Range("C3").Value = 30
$ws1.Range("A4").Value = "North Ame"
$ws1.Range("B4").Value = "Flamie"
$ws1.Range("A4").Select()

# Sheet "g2"
$ws2 = $wb.Worksheets.Item("g2")
$ws2.Range("C2").Value = 40
$ws2.Range("C3").Value = 40
$ws2.Range("A4").Value = "North Ame"
$ws2.Range("B4").Value = "Flamie"
$ws2.Range("A4").Select()

# Sheet "faze"
$ws3 = $wb.Worksheets.Item("faze")
$ws3.Range("A4").Value = "North Ame"
$ws3.Range("B4").Value = "Flamie"
$ws3.Range("A4").Select()
